$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers
$ws.Range("F1").Value = "height"
$ws.Range("G1").Value = "weight"

# Copy style from E1 (existing header style) to F1:G1
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Update values for rows 2-7
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667   # column E: fantasy points
    $ws.Cells.Item($r, 6).Value = 252                  # column F: height
    $ws.Cells.Item($r, 7).Value = 0                    # column G: weight
}
